$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feedback")
$ws.Activate()

# Add "board" feedback markers ("X") in the new H column for the rows
# that correspond to the "board" related events.
$ws.Range("H8").Value = "X"
$ws.Range("H10").Value = "X"
$ws.Range("H11").Value = "X"
$ws.Range("H12").Value = "X"

# Leave selection on the last edited cell, as in the authored change.
$ws.Range("H12").Select()
